$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("C22").Select()
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "To do"
$newSheet.Activate()
Write-Output "done"
